$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-28 Sunday", "2025-12-29 Monday"),
    @("23×33=", "88×93="),
    @("82×90=", "42×57="),
    @("72×21=", "32×93="),
    @("96×79=", "13×52="),
    @("82×85=", "37×37="),
    @("63×85=", "73×27="),
    @("58×57=", "22×58="),
    @("65×20=", "83×81="),
    @("54×14=", "90×85="),
    @("77×62=", "40×15="),
    @("11×55=", "26×87="),
    @("87×35=", "43×16="),
    @("19×71=", "41×19="),
    @("28×46=", "28×91="),
    @("44×25=", "52×95="),
    @("94×26=", "30×11="),
    @("38×68=", "93×95="),
    @("72×23=", "29×51="),
    @("35×68=", "60×68="),
    @("34×48=", "30×54="),
    @("19×50=", "36×35="),
    @("83×66=", "24×51="),
    @("99×33=", "82×46="),
    @("71×73=", "23×27="),
    @("61×97=", "14×36=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
